$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- Sheet ALC ---
# Row 13 (item id 2144)
$ws1.Range("H13").Value = 70006
$ws1.Range("I13").Value = 0
$ws1.Range("K13").Value = 0
$ws1.Range("M13").ClearContents()

# Row 87 (item id 10651)
$ws1.Range("H87").Value = 34666.668
$ws1.Range("I87").Value = 0
$ws1.Range("J87").Value = 34666.668
$ws1.Range("K87").Value = 0
$ws1.Range("L87").Value = 34666.668
$ws1.Range("N87").Value = -37162.668
$ws1.Range("M87").ClearContents()

# Row 90 (item id 10651)
$ws1.Range("H90").Value = 34666.668
$ws1.Range("I90").Value = 0
$ws1.Range("J90").Value = 34666.668
$ws1.Range("K90").Value = 0
$ws1.Range("L90").Value = 104000.004
$ws1.Range("N90").Value = -116480.004
$ws1.Range("M90").ClearContents()

# Row 96 (item id 19894)
$ws1.Range("H96").Value = 1069210.2
$ws1.Range("I96").Value = 1972766.8
$ws1.Range("J96").Value = 1370.8182
$ws1.Range("K96").Value = 5918300.4
$ws1.Range("L96").Value = 4112.4546
$ws1.Range("M96").Value = -5916927.4
$ws1.Range("N96").Value = -6858.4546

# Row 98 (item id 36237)
$ws1.Range("H98").Value = 1745.1
$ws1.Range("I98").Value = 1645.4546
$ws1.Range("J98").Value = 1866.8889
$ws1.Range("K98").Value = 1645.4546
$ws1.Range("L98").Value = 1866.8889
$ws1.Range("M98").Value = -147.4546
$ws1.Range("N98").Value = -4862.8889

# Row 122 (item id 36237)
$ws1.Range("H122").Value = 1745.1
$ws1.Range("I122").Value = 1645.4546
$ws1.Range("J122").Value = 1866.8889
$ws1.Range("K122").Value = 4936.3638
$ws1.Range("L122").Value = 5600.6667
$ws1.Range("M122").Value = -2486.3638
$ws1.Range("N122").Value = -10500.6667

# Row 132 (item id 44049)
$ws1.Range("H132").Value = 6491.8
$ws1.Range("I132").Value = 4989.778
$ws1.Range("J132").Value = 12499.889
$ws1.Range("K132").Value = 14969.334
$ws1.Range("L132").Value = 37499.667
$ws1.Range("M132").Value = -12439.334
$ws1.Range("N132").Value = -42559.667

# --- Sheet ARM ---
# Row 45 (item id 27714)
$ws2.Range("H45").Value = 2566.8
$ws2.Range("I45").Value = 1360.6471
$ws2.Range("J45").Value = 4144.077
$ws2.Range("K45").Value = 1360.6471
$ws2.Range("L45").Value = 4144.077
$ws2.Range("M45").Value = -983.6470999999999
$ws2.Range("N45").Value = -4898.077

# Row 81 (item id 10841)
$ws2.Range("H81").Value = 0
$ws2.Range("J81").Value = 0
$ws2.Range("L81").Value = 0
$ws2.Range("N81").ClearContents()

# Row 84 (item id 10841)
$ws2.Range("H84").Value = 0
$ws2.Range("J84").Value = 0
$ws2.Range("L84").Value = 0
$ws2.Range("N84").ClearContents()

# Row 122 (item id 36168)
$ws2.Range("H122").Value = 2729.2778
$ws2.Range("I122").Value = 1416.2858
$ws2.Range("J122").Value = 3564.818
$ws2.Range("K122").Value = 4248.857400000001
$ws2.Range("L122").Value = 10694.454
$ws2.Range("M122").Value = -1798.857400000001
$ws2.Range("N122").Value = -15594.454

# Row 132 (item id 43997)
$ws2.Range("H132").Value = 2122224.5
$ws2.Range("I132").Value = 1532.9429
$ws2.Range("K132").Value = 4598.8287
$ws2.Range("M132").Value = -2068.8287

# --- Sheet BSM ---
# Row 20 (item id 14149)
$ws3.Range("H20").Value = 4040.625
$ws3.Range("I20").Value = 5904
$ws3.Range("J20").Value = 3419.5
$ws3.Range("K20").Value = 5904
$ws3.Range("L20").Value = 3419.5
$ws3.Range("M20").Value = -5657
$ws3.Range("N20").Value = -3913.5

# Row 134 (item id 43998)
$ws3.Range("H134").Value = 5057.298
$ws3.Range("I134").Value = 2204.6365
$ws3.Range("J134").Value = 7567.64
$ws3.Range("K134").Value = 6613.9095
$ws3.Range("L134").Value = 22702.92
$ws3.Range("M134").Value = -4078.9095
$ws3.Range("N134").Value = -27772.92

# --- Sheet CRP ---
# Row 16 (item id 27691)
$ws4.Range("H16").Value = 3228.2856
$ws4.Range("I16").Value = 4520.1665
$ws4.Range("K16").Value = 4520.1665
$ws4.Range("M16").Value = -4233.1665

# Row 86 (item id 12584)
$ws4.Range("H86").Value = 34486736
$ws4.Range("I86").Value = 66669280
$ws4.Range("J86").Value = 5443.4287
$ws4.Range("K86").Value = 66669280
$ws4.Range("L86").Value = 5443.4287
$ws4.Range("M86").Value = -66668157
$ws4.Range("N86").Value = -7689.4287

# Row 89 (item id 12584)
$ws4.Range("H89").Value = 34486736
$ws4.Range("I89").Value = 66669280
$ws4.Range("J89").Value = 5443.4287
$ws4.Range("K89").Value = 333346400
$ws4.Range("L89").Value = 27217.1435
$ws4.Range("M89").Value = -333340784
$ws4.Range("N89").Value = -38449.14350000001

# Row 113 (item id 27691)
$ws4.Range("H113").Value = 3228.2856
$ws4.Range("I113").Value = 4520.1665
$ws4.Range("K113").Value = 4520.1665
$ws4.Range("M113").Value = -2350.1665

# --- Sheet CUL ---
# Row 118 (item id 27872)
$ws5.Range("H118").Value = 5490
$ws5.Range("I118").Value = 0
$ws5.Range("K118").Value = 0
$ws5.Range("M118").ClearContents()

# --- Sheet GSM ---
# Row 92 (item id 18094)
$ws6.Range("H92").Value = 0
$ws6.Range("J92").Value = 0
$ws6.Range("L92").Value = 0
$ws6.Range("N92").ClearContents()

# Row 102 (item id 36169)
$ws6.Range("H102").Value = 1744254.2
$ws6.Range("I102").Value = 2552569.2
$ws6.Range("J102").Value = 3267.7693
$ws6.Range("K102").Value = 2552569.2
$ws6.Range("L102").Value = 3267.7693
$ws6.Range("M102").Value = -2550947.2
$ws6.Range("N102").Value = -6511.7693

# Row 122 (item id 36182)
$ws6.Range("H122").Value = 3767.5
$ws6.Range("I122").Value = 3027.7778
$ws6.Range("J122").Value = 5986.6665
$ws6.Range("K122").Value = 9083.3334
$ws6.Range("L122").Value = 17959.9995
$ws6.Range("M122").Value = -6633.3334
$ws6.Range("N122").Value = -22859.9995

# Row 126 (item id 36184)
$ws6.Range("H126").Value = 50006400
$ws6.Range("J126").Value = 8000
$ws6.Range("L126").Value = 24000
$ws6.Range("N126").Value = -28940

# Row 132 (item id 44008)
$ws6.Range("H132").Value = 869760.5600000001
$ws6.Range("I132").Value = 1345511.6
$ws6.Range("K132").Value = 4036534.8
$ws6.Range("M132").Value = -4034004.8

# --- Sheet LTW ---
# Row 7 (item id 36249)
$ws7.Range("H7").Value = 1625.4286
$ws7.Range("I7").Value = 1627.591
$ws7.Range("J7").Value = 1617.5
$ws7.Range("K7").Value = 1627.591
$ws7.Range("L7").Value = 1617.5
$ws7.Range("M7").Value = -1515.591
$ws7.Range("N7").Value = -1841.5

# Row 16 (item id 5289)
$ws7.Range("H16").Value = 8800.786
$ws7.Range("I16").Value = 1785.4615
$ws7.Range("J16").Value = 100000
$ws7.Range("K16").Value = 1785.4615
$ws7.Range("L16").Value = 100000
$ws7.Range("M16").Value = -1615.4615
$ws7.Range("N16").Value = -100340

# Row 36 (item id 34261)
$ws7.Range("H36").Value = 30000
$ws7.Range("J36").Value = 30000
$ws7.Range("L36").Value = 30000
$ws7.Range("N36").Value = -31124

# Row 61 (item id 27740)
$ws7.Range("H61").Value = 5722.5
$ws7.Range("I61").Value = 1926.6666
$ws7.Range("J61").Value = 8000
$ws7.Range("K61").Value = 1926.6666
$ws7.Range("L61").Value = 8000
$ws7.Range("M61").Value = -1724.6666
$ws7.Range("N61").Value = -8404

# Row 113 (item id 27740)
$ws7.Range("H113").Value = 5722.5
$ws7.Range("I113").Value = 1926.6666
$ws7.Range("J113").Value = 8000
$ws7.Range("K113").Value = 1926.6666
$ws7.Range("L113").Value = 8000
$ws7.Range("M113").Value = 243.3334
$ws7.Range("N113").Value = -12340

# Row 122 (item id 36247)
$ws7.Range("H122").Value = 22279.857
$ws7.Range("I122").Value = 29124.75
$ws7.Range("J122").Value = 13153.333
$ws7.Range("K122").Value = 87374.25
$ws7.Range("L122").Value = 39459.999
$ws7.Range("M122").Value = -84924.25
$ws7.Range("N122").Value = -44359.999

# Row 126 (item id 36249)
$ws7.Range("H126").Value = 1625.4286
$ws7.Range("I126").Value = 1627.591
$ws7.Range("J126").Value = 1617.5
$ws7.Range("K126").Value = 4882.772999999999
$ws7.Range("L126").Value = 4852.5
$ws7.Range("M126").Value = -2412.772999999999
$ws7.Range("N126").Value = -9792.5

# --- Sheet WVR ---
# Row 51 (item id 3162)
$ws8.Range("H51").Value = 0
$ws8.Range("I51").Value = 0
$ws8.Range("K51").Value = 0
$ws8.Range("M51").ClearContents()

# Row 132 (item id 44029)
$ws8.Range("H132").Value = 2291.7856
$ws8.Range("I132").Value = 1574.12
$ws8.Range("J132").Value = 3347.1765
$ws8.Range("K132").Value = 4722.36
$ws8.Range("L132").Value = 10041.5295
$ws8.Range("M132").Value = -2192.36
$ws8.Range("N132").Value = -15101.5295

# Row 136 (item id 44031)
$ws8.Range("H136").Value = 9626547
$ws8.Range("I136").Value = 14721728
$ws8.Range("J136").Value = 2315
$ws8.Range("K136").Value = 44165184
$ws8.Range("L136").Value = 6945
$ws8.Range("M136").Value = -44162634
$ws8.Range("N136").Value = -12045
